$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("A18").Value = "billion 2023 dollars"
$ws.Range("A21").Value = "million 2023 dollars"
$ws.Range("A24").Value = "2023 dollars"
$ws.Range("B29").Value = 'which in this case is "2012 dollars per 2023 dollar."'

[void]$ws.Range("B30").Select()
